$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
